$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-111). The diff bumps each of these values from 45171 to 45172
# (i.e. one day later), leaving everything else untouched.
for ($row = 2; $row -le 111; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
